$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.321.60"
$ws.Range("E2").Value = "'  -0.05%  "
$ws.Range("D3").Value = "'1.933.83"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("D4").Value = "'0.9965"
$ws.Range("E4").Value = "'  -0.42%  "
$ws.Range("D5").Value = "'0.7684"
$ws.Range("E5").Value = "'  +6.01%  "
$ws.Range("D6").Value = "'248.43"
$ws.Range("E6").Value = "'  -1.25%  "
$ws.Range("D7").Value = "'0.9969"
$ws.Range("E7").Value = "'  -0.40%  "
$ws.Range("D8").Value = "'28.41"
$ws.Range("E8").Value = "'  +0.65%  "
$ws.Range("D9").Value = "'0.3231"
$ws.Range("E9").Value = "'  -2.85%  "
$ws.Range("D10").Value = "'0.07124"
$ws.Range("E10").Value = "'  -1.47%  "
$ws.Range("D11").Value = "'0.7929"
$ws.Range("E11").Value = "'  -2.56%  "
$ws.Range("E12").Value = "'  -1.29%  "
$ws.Range("D13").Value = "'1.933.65"
$ws.Range("E13").Value = "'  -0.23%  "
$ws.Range("D14").Value = "'5.385"
$ws.Range("E14").Value = "'  -1.97%  "
$ws.Range("D15").Value = "'94.84"
$ws.Range("E15").Value = "'  +0.04%  "
$ws.Range("D16").Value = "'14.77"
$ws.Range("E16").Value = "'  -2.90%  "
$ws.Range("D17").Value = "'30.318.54"
$ws.Range("E17").Value = "'  -0.11%  "
$ws.Range("D18").Value = "'253.52"
$ws.Range("E18").Value = "'  +1.27%  "
$ws.Range("D19").Value = "'0.000008048"
$ws.Range("E19").Value = "'  -3.04%  "
$ws.Range("D20").Value = "'5.815"
$ws.Range("D21").Value = "'2.186.82"
$ws.Range("E21").Value = "'  -0.18%  "
$ws.Range("D22").Value = "'0.9966"
$ws.Range("E22").Value = "'  -0.41%  "
$ws.Range("D23").Value = "'0.9955"
$ws.Range("E23").Value = "'  -0.46%  "
$ws.Range("D24").Value = "'6.850"
$ws.Range("E24").Value = "'  -2.16%  "
$ws.Range("D25").Value = "'9.601"
$ws.Range("E25").Value = "'  -1.77%  "
$ws.Range("D26").Value = "'165.17"
$ws.Range("E26").Value = "'  +0.75%  "
$ws.Range("D27").Value = "'0.1366"
$ws.Range("E27").Value = "'  +3.14%  "
$ws.Range("E28").Value = "'  -2.72%  "
$ws.Range("D29").Value = "'19.12"
$ws.Range("E29").Value = "'  -1.13%  "
$ws.Range("D30").Value = "'1.377"
$ws.Range("E30").Value = "'  +1.89%  "
$ws.Range("D31").Value = "'1.528"
$ws.Range("E31").Value = "'  -2.71%  "
$ws.Range("D32").Value = "'4.439"
$ws.Range("E32").Value = "'  -0.12%  "
$ws.Range("D33").Value = "'4.152"
$ws.Range("E33").Value = "'  -0.67%  "
$ws.Range("D34").Value = "'0.05180"
$ws.Range("E34").Value = "'  -0.63%  "
$ws.Range("D35").Value = "'1.297"
$ws.Range("E35").Value = "'  +0.48%  "
$ws.Range("D36").Value = "'0.7538"
$ws.Range("E36").Value = "'  +0.35%  "
$ws.Range("D37").Value = "'2.763"
$ws.Range("E37").Value = "'  +0.55%  "
$ws.Range("D38").Value = "'0.01967"
$ws.Range("E38").Value = "'  -0.87%  "
$ws.Range("D39").Value = "'2.806"
$ws.Range("E39").Value = "'  -0.93%  "
$ws.Range("D40").Value = "'78.45"
$ws.Range("E40").Value = "'  -2.77%  "
$ws.Range("D41").Value = "'6.426"
$ws.Range("E41").Value = "'  -0.31%  "
$ws.Range("D42").Value = "'0.4538"
$ws.Range("E42").Value = "'  -0.27%  "
$ws.Range("D43").Value = "'2.000"
$ws.Range("E43").Value = "'  -2.06%  "
$ws.Range("D44").Value = "'0.9975"
$ws.Range("E44").Value = "'  -0.28%  "
$ws.Range("D45").Value = "'0.8357"
$ws.Range("E45").Value = "'  -1.35%  "
$ws.Range("D46").Value = "'102.65"
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.565"
$ws.Range("E47").Value = "'  +1.23%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.826"
$ws.Range("E48").Value = "'  +0.16%  "
$ws.Range("D49").Value = "'992.78"
$ws.Range("E49").Value = "'  +13.29%  "
$ws.Range("E50").Value = "'  +1.53%  "
$ws.Range("B51").Value = "'Decentraland"
$ws.Range("C51").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4191"
$ws.Range("E51").Value = "'  -0.05%  "
